$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update F column "want to go" counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 30
$ws1.Range("F3").Value = 49
$ws1.Range("F5").Value = 182
$ws1.Range("F6").Value = 1076
$ws1.Range("F7").Value = 1045
$ws1.Range("F8").Value = 8172
$ws1.Range("F10").Value = 208
$ws1.Range("F11").Value = 6896
$ws1.Range("F12").Value = 171
$ws1.Range("F13").Value = 300
$ws1.Range("F14").Value = 4994
$ws1.Range("F16").Value = 5432
$ws1.Range("F17").Value = 1072
$ws1.Range("F18").Value = 330
$ws1.Range("F19").Value = 337
$ws1.Range("F20").Value = 460
$ws1.Range("F22").Value = 254
$ws1.Range("F26").Value = 9167
$ws1.Range("F27").Value = 71
$ws1.Range("F28").Value = 1664
$ws1.Range("F29").Value = 764
$ws1.Range("F30").Value = 40
$ws1.Range("F32").Value = 1859
$ws1.Range("F33").Value = 71
$ws1.Range("F34").Value = 78
$ws1.Range("F36").Value = 1009
$ws1.Range("F37").Value = 1877
$ws1.Range("F40").Value = 4781
$ws1.Range("F42").Value = 1161
$ws1.Range("F43").Value = 76
$ws1.Range("F46").Value = 36
$ws1.Range("F47").Value = 914
$ws1.Range("F48").Value = 1252

# Sheet "全部类型" (All types) - update F column "want to go" counts
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 30
$ws4.Range("F4").Value = 49
$ws4.Range("F6").Value = 182
$ws4.Range("F8").Value = 1076
$ws4.Range("F9").Value = 1045
$ws4.Range("F10").Value = 8172
$ws4.Range("F12").Value = 208
$ws4.Range("F13").Value = 6896
$ws4.Range("F14").Value = 171
$ws4.Range("F15").Value = 300
$ws4.Range("F17").Value = 4994
$ws4.Range("F19").Value = 5432
$ws4.Range("F20").Value = 1072
$ws4.Range("F21").Value = 330
$ws4.Range("F22").Value = 337
$ws4.Range("F23").Value = 460
$ws4.Range("F25").Value = 254
$ws4.Range("F27").Value = 9167
$ws4.Range("F28").Value = 71
$ws4.Range("F29").Value = 1664
$ws4.Range("F30").Value = 764
$ws4.Range("F31").Value = 40
$ws4.Range("F33").Value = 1859
$ws4.Range("F34").Value = 71
$ws4.Range("F35").Value = 78
$ws4.Range("F37").Value = 1009
$ws4.Range("F38").Value = 1877
$ws4.Range("F41").Value = 4782
$ws4.Range("F43").Value = 1161
$ws4.Range("F44").Value = 76
$ws4.Range("F47").Value = 914
$ws4.Range("F48").Value = 1252
